$d = $word.ActiveDocument

$replacements = @(
    @("954÷7=", "216÷8="),
    @("946÷8=", "740÷8="),
    @("703÷9=", "681÷8="),
    @("526÷5=", "484÷9="),
    @("517÷2=", "858÷3="),
    @("648÷3=", "342÷8="),
    @("307÷7=", "620÷3="),
    @("512÷2=", "568÷4="),
    @("275÷9=", "780÷8="),
    @("438÷3=", "581÷8="),
    @("529÷3=", "469÷2="),
    @("280÷6=", "897÷8="),
    @("444÷6=", "104÷6="),
    @("584÷4=", "674÷6="),
    @("719÷8=", "758÷3="),
    @("690÷4=", "259÷2="),
    @("636÷2=", "809÷6="),
    @("113÷8=", "373÷2="),
    @("368÷8=", "486÷3="),
    @("631÷5=", "823÷2="),
    @("433÷6=", "695÷8="),
    @("563÷6=", "575÷7="),
    @("430÷6=", "677÷8="),
    @("257÷9=", "366÷7="),
    @("638÷9=", "898÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
